$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Institute")
$ws.Range("A1").Value = "test"
